$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Simple text replacements (in-place within existing runs/paragraphs)
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "Senior Software Engineer with 21 years of experience",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Senior Software Engineer with 15+ years of experience",
    2)

$d.Content.Find.Execute(
    "Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys",
    2)

$d.Content.Find.Execute(
    "Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands of users simultaneously",
    2)

$d.Content.Find.Execute(
    "Integrated mapping and visualization tools for political campaign data analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs",
    2)

# ---------------------------------------------------------------------------
# 2) Structural changes (paragraph insertions / removals).
#    Performed from the bottom of the document upward so that earlier
#    paragraph indices remain valid while later ones are being edited.
# ---------------------------------------------------------------------------

# --- Replace the EDUCATION section (Heading2 "EDUCATION" + two Heading3
#     degree lines) with a single plain bullet paragraph. -------------------
$eduIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.Trim() -eq "EDUCATION") {
        $eduIndex = $i
        break
    }
}
$d.Paragraphs.Item($eduIndex + 2).Range.Delete()
$d.Paragraphs.Item($eduIndex + 1).Range.Delete()
$d.Paragraphs.Item($eduIndex).Range.Delete()
$d.Paragraphs.Item($eduIndex - 1).Range.InsertParagraphAfter()
$d.Paragraphs.Item($eduIndex).Range.Text = "• Trained staff on PHP/MySQL for data analysis and reporting systems"

# --- Add a new bullet after "Developed innovative approaches to
#     visualizing demographic and market data..." -------------------------
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.Trim() -eq "• Developed innovative approaches to visualizing demographic and market data for enhanced client understanding") {
        $anchorIndex = $i
        break
    }
}
$d.Paragraphs.Item($anchorIndex).Range.InsertParagraphAfter()
$d.Paragraphs.Item($anchorIndex + 1).Range.Text = "• Trained staff on building Python tooling for report generation and analysis"

# --- Add a new bullet after "Managed technology infrastructure supporting
#     community health initiatives across multiple countries" -------------
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.Trim() -eq "• Managed technology infrastructure supporting community health initiatives across multiple countries") {
        $anchorIndex = $i
        break
    }
}
$d.Paragraphs.Item($anchorIndex).Range.InsertParagraphAfter()
$d.Paragraphs.Item($anchorIndex + 1).Range.Text = "• Architected and developed 25 Drupal sites to integrate with membership databases, activism CRMs and government agencies, under guidelines from Kellogg Foundation and Robert Wood Johnson Foundation"

# --- Add a new bullet after "Collaborated with political strategists to
#     translate geospatial requirements into technical solutions" ---------
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.Trim() -eq "• Collaborated with political strategists to translate geospatial requirements into technical solutions") {
        $anchorIndex = $i
        break
    }
}
$d.Paragraphs.Item($anchorIndex).Range.InsertParagraphAfter()
$d.Paragraphs.Item($anchorIndex + 1).Range.Text = "• Handled billions of records with millions of columns in high-performance CRM system"
